# Strengthen the EMM: bump the "High severity" relative-risk parameters
# on the SimParameters sheet from 0.75 to 0.8, for both the Abortion
# and Preeclampsia blocks. Downstream formulas on the potential_preg_trt
# and potential_preec_trt sheets reference these cells and will
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SimParameters")

$ws.Range("B6").Value = 0.8
$ws.Range("B10").Value = 0.8

# Match the author's final cell selection on the SimParameters sheet.
$ws.Activate()
$ws.Range("B11").Select()

$excel.CalculateFullRebuild()
